# Insert a new data row at row 726 (pushing existing rows 726-837 down to 727-838)
# and populate it with the new "Cebollín" price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(726).Insert()

$ws.Cells.Item(726, 1).Value2 = 3
$ws.Cells.Item(726, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(726, 3).Value2 = "Coquimbo"
$ws.Cells.Item(726, 4).Value2 = 45218
$ws.Cells.Item(726, 5).Value2 = 5
$ws.Cells.Item(726, 6).Value2 = 100112037
$ws.Cells.Item(726, 7).Value2 = "Cebollín"
$ws.Cells.Item(726, 8).Value2 = "Sin especificar"
$ws.Cells.Item(726, 9).Value2 = "Primera"
$ws.Cells.Item(726, 10).Value2 = 193
$ws.Cells.Item(726, 11).Value2 = 4000
$ws.Cells.Item(726, 12).Value2 = 4500
$ws.Cells.Item(726, 13).Value2 = 4202
$ws.Cells.Item(726, 14).Value2 = "`$/paquete 36 unidades"
$ws.Cells.Item(726, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(726, 16).Value2 = 117
$ws.Cells.Item(726, 17).Value2 = 36
$ws.Cells.Item(726, 18).Value2 = "Hortaliza"
